$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C5").ClearContents()
$ws.Range("B6").Value = 5000

# Copy A1 (which has border1 = full box) style into C5 via PasteSpecial formats
$ws.Range("A1").Copy()
$ws.Range("C5").PasteSpecial(-4122)  # xlPasteFormats = -4122

# Now clear the border specifically (all edges none)
$ws.Range("C5").Borders.LineStyle = -4142

Write-Host "C5 border left after clear:" $ws.Range("C5").Borders.Item(7).LineStyle
